$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new header cells for the added "Tag" and "Instrument" columns
$ws.Range("K1").Value = "Tag"
$ws.Range("L1").Value = "Instrument"

# Match the style of other non-header cells (plain Arial 10pt, black, no bold/border)
$ws.Range("K1:L1").Font.Name = "Arial"
$ws.Range("K1:L1").Font.Size = 10
$ws.Range("K1:L1").Font.Bold = $false
$ws.Range("K1:L1").Font.Color = 0

# Update the selection to match the target workbook
$ws.Range("K2").Select()
